$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-11 from 45224 to 45233
$ws.Range("C2:C11").Value = 45233
